$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "https://support.google.com/docs/answer/3093281?hl=en"
$ws.Range("B4").Value = "GOOGLE FINANCE"

$ws.Range("B4").Select()
